$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.130.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.230.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.09"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.594"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.234.36"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.604"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.07"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.133"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.11"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.750.90"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.234.76"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.038.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.15"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.965"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.61"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.16"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.02"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.13"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.43"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "638.48"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.43"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.21"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.90"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.63"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.375"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.60"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.122"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.884.46"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.95%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0392"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.95%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.06"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.08%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.92%  "
